$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings with rich-text runs) ---
# "Volume 30   Number  9" -> "Volume 30   Number  10"
$ws.Range("A8").Value = "Volume 30   Number  10"

# "Report Covering the Week  2/27/2023  Through  3/5/2023" -> "...3/6/2023  Through  3/12/2023"
$ws.Range("C9").Value = "Report Covering the Week  3/6/2023  Through  3/12/2023"

# --- Row 14 ---
$ws.Range("L14").Value = -60

# --- Row 15 ---
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = "***.*"
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 6
$ws.Range("K15").Value = 20
$ws.Range("L15").Value = -40
$ws.Range("M15").Value = -40
$ws.Range("N15").Value = -71.428571428571

# --- Row 16 ---
$ws.Range("F16").Value = 23
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = 43.75
$ws.Range("I16").Value = 71
$ws.Range("J16").Value = 41
$ws.Range("K16").Value = 73.170731707317
$ws.Range("L16").Value = 121.875
$ws.Range("M16").Value = 10.9375
$ws.Range("N16").Value = -76.872964169381

# --- Row 17 ---
$ws.Range("C17").Value = 19
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 280
$ws.Range("F17").Value = 59
$ws.Range("G17").Value = 40
$ws.Range("H17").Value = 47.5
$ws.Range("I17").Value = 163
$ws.Range("J17").Value = 119
$ws.Range("K17").Value = 36.974789915966
$ws.Range("L17").Value = 59.803921568627
$ws.Range("M17").Value = 64.646464646464
$ws.Range("N17").Value = -18.090452261306

# --- Row 18 ---
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 9
$ws.Range("E18").Value = -44.444444444444
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = 8.333333333333
$ws.Range("I18").Value = 69
$ws.Range("J18").Value = 43
$ws.Range("K18").Value = 60.465116279069
$ws.Range("L18").Value = 27.777777777777
$ws.Range("M18").Value = -28.125
$ws.Range("N18").Value = -89.351851851851

# --- Row 19 ---
$ws.Range("C19").Value = 27
$ws.Range("D19").Value = 32
$ws.Range("E19").Value = -15.625
$ws.Range("F19").Value = 108
$ws.Range("G19").Value = 109
$ws.Range("H19").Value = -0.917431192660
$ws.Range("I19").Value = 285
$ws.Range("J19").Value = 285
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 50
$ws.Range("M19").Value = 68.639053254437
$ws.Range("N19").Value = -1.724137931034

# --- Row 20 ---
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 15
$ws.Range("E20").Value = -53.333333333333
$ws.Range("F20").Value = 29
$ws.Range("G20").Value = 36
$ws.Range("H20").Value = -19.444444444444
$ws.Range("I20").Value = 73
$ws.Range("J20").Value = 81
$ws.Range("K20").Value = -9.876543209876
$ws.Range("L20").Value = 82.5
$ws.Range("M20").Value = 23.728813559322
$ws.Range("N20").Value = -92.603850050658

# --- Row 21 ---
$ws.Range("C21").Value = 62
$ws.Range("D21").Value = 65
$ws.Range("E21").Value = -4.615384615384
$ws.Range("F21").Value = 249
$ws.Range("G21").Value = 226
$ws.Range("H21").Value = 10.176991150442
$ws.Range("I21").Value = 669
$ws.Range("J21").Value = 575
$ws.Range("K21").Value = 16.347826086956
$ws.Range("L21").Value = 54.503464203233
$ws.Range("M21").Value = 34.068136272545
$ws.Range("N21").Value = -72.771672771672

# --- Row 23 ---
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 200
$ws.Range("F23").Value = 9
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 28.571428571428
$ws.Range("I23").Value = 33
$ws.Range("J23").Value = 19
$ws.Range("K23").Value = 73.684210526315
$ws.Range("L23").Value = 94.117647058823
$ws.Range("M23").Value = 312.5

# --- Row 24 ---
$ws.Range("C24").Value = 81
$ws.Range("D24").Value = 67
$ws.Range("E24").Value = 20.895522388059
$ws.Range("F24").Value = 289
$ws.Range("G24").Value = 271
$ws.Range("H24").Value = 6.642066420664
$ws.Range("I24").Value = 766
$ws.Range("J24").Value = 619
$ws.Range("K24").Value = 23.747980613893
$ws.Range("L24").Value = 69.094922737306
$ws.Range("M24").Value = 21.780604133545

# --- Row 25 ---
$ws.Range("C25").Value = 33
$ws.Range("D25").Value = 30
$ws.Range("E25").Value = 10
$ws.Range("G25").Value = 134
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 322
$ws.Range("J25").Value = 291
$ws.Range("K25").Value = 10.652920962199
$ws.Range("L25").Value = 84
$ws.Range("M25").Value = -3.012048192771

# --- Row 26 ---
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = -33.333333333333
$ws.Range("F26").Value = 5
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 11
$ws.Range("J26").Value = 14
$ws.Range("K26").Value = -21.428571428571
$ws.Range("L26").Value = -31.25

# --- Row 27 ---
$ws.Range("C27").Value = 7
$ws.Range("E27").Value = 133.333333333333
$ws.Range("F27").Value = 14
$ws.Range("G27").Value = 10
$ws.Range("H27").Value = 40
$ws.Range("I27").Value = 34
$ws.Range("J27").Value = 29
$ws.Range("K27").Value = 17.241379310344
$ws.Range("L27").Value = 47.826086956521

# --- Row 28 ---
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = "***.*"
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = 100
$ws.Range("L28").Value = -40
$ws.Range("N28").Value = -83.333333333333

# --- Row 29 ---
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = "***.*"
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = 100
$ws.Range("L29").Value = -40
$ws.Range("N29").Value = -81.25

# --- Row 30 ---
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = -100
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = -100
$ws.Range("J30").Value = 3
$ws.Range("K30").Value = -66.666666666666
